$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9542087912559509
$ws.Range("B1").Value = 0.9135027527809143
$ws.Range("C1").Value = 0.9288582801818848
$ws.Range("D1").Value = 1.155512094497681
$ws.Range("E1").Value = 1.045491456985474
